# Updated cryptos list on Sat Oct  7 10:36:24 UTC 2023 with GitHub Actions
# Applies latest price / 1h-volume snapshot (and two coin-row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Row, $Col, $Text)
    $c = $ws.Cells.Item($Row, $Col)
    $c.NumberFormat = "@"
    $c.Value = $Text
    $c.ClearFormats()
}

# Row 2 - Bitcoin
Set-CellText 2 4 "27.955.09"
Set-CellText 2 5 "  +0.69%  "

# Row 3 - Ethereum
Set-CellText 3 4 "1.640.54"
Set-CellText 3 5 "  +0.28%  "

# Row 4 - TetherUSD
Set-CellText 4 5 "  +0.14%  "

# Row 5 - BNB
Set-CellText 5 4 "213.07"
Set-CellText 5 5 "  +0.41%  "

# Row 6 - XRP
Set-CellText 6 4 "0.523"
Set-CellText 6 5 "  +0.05%  "

# Row 7 - USDC
Set-CellText 7 5 "  +0.13%  "

# Row 8 - Solana
Set-CellText 8 4 "23.57"
Set-CellText 8 5 "  +0.46%  "

# Row 9 - Cardano
Set-CellText 9 5 "  -1.47%  "

# Row 10 - Dogecoin
Set-CellText 10 5 "  +0.35%  "

# Row 11 - TRON
Set-CellText 11 5 "  +2.44%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-CellText 12 5 "  +0.27%  "

# Row 13 - WrappedEther
Set-CellText 13 4 "1.641.77"
Set-CellText 13 5 "  +0.25%  "

# Row 14 - Polygon->Polkadot
Set-CellText 14 2 "Polkadot"
Set-CellText 14 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText 14 4 "4.09"
Set-CellText 14 5 "  +1.10%  "

# Row 15 - Polkadot->Polygon
Set-CellText 15 2 "Polygon"
Set-CellText 15 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-CellText 15 4 "0.575"
Set-CellText 15 5 "  +3.54%  "

# Row 16 - Litecoin
Set-CellText 16 4 "65.88"
Set-CellText 16 5 "  +1.08%  "

# Row 17 - WrappedBTC
Set-CellText 17 4 "27.960.54"
Set-CellText 17 5 "  +0.88%  "

# Row 18 - BitcoinCash
Set-CellText 18 4 "233.28"
Set-CellText 18 5 "  +0.84%  "

# Row 19 - ShibaInu
Set-CellText 19 4 "0.0₃0724"
Set-CellText 19 5 "  +0.51%  "

# Row 20 - Chainlink
Set-CellText 20 4 "7.61"
Set-CellText 20 5 "  +0.01%  "

# Row 21 - Dai
Set-CellText 21 5 "  +0.12%  "

# Row 22 - Avalanche
Set-CellText 22 4 "10.71"
Set-CellText 22 5 "  +0.00%  "

# Row 23 - Uniswap
Set-CellText 23 5 "  -0.07%  "

# Row 24 - Toncoin
Set-CellText 24 5 "  -1.81%  "

# Row 25 - Monero
Set-CellText 25 4 "151.22"
Set-CellText 25 5 "  +1.42%  "

# Row 26 - Cosmos
Set-CellText 26 5 "  +1.38%  "

# Row 27 - EthereumClassic
Set-CellText 27 4 "15.69"
Set-CellText 27 5 "  +0.65%  "

# Row 28 - Stellar
Set-CellText 28 5 "  -0.25%  "

# Row 29 - BinanceUSD
Set-CellText 29 5 "  +0.15%  "

# Row 30 - PancakeSwap
Set-CellText 30 5 "  +0.49%  "

# Row 31 - Hedera
Set-CellText 31 5 "  -0.08%  "

# Row 32 - Filecoin
Set-CellText 32 5 "  +1.69%  "

# Row 33 - InternetComputer(DFINITY)
Set-CellText 33 5 "  +0.85%  "

# Row 34 - Maker
Set-CellText 34 4 "1.407.74"
Set-CellText 34 5 "  -5.08%  "

# Row 35 - LidoDAOToken
Set-CellText 35 4 "1.57"
Set-CellText 35 5 "  +1.54%  "

# Row 36 - HuobiToken
Set-CellText 36 5 "  +1.38%  "

# Row 37 - VeChain->ARBITRUM
Set-CellText 37 2 "ARBITRUM"
Set-CellText 37 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText 37 4 "0.883"
Set-CellText 37 5 "  +0.01%  "

# Row 38 - ARBITRUM->VeChain
Set-CellText 38 2 "VeChain"
Set-CellText 38 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText 38 4 "0.0169"
Set-CellText 38 5 "  +0.94%  "

# Row 39 - ImmutableX
Set-CellText 39 4 "0.557"
Set-CellText 39 5 "  -0.37%  "

# Row 40 - TrustWalletToken
Set-CellText 40 4 "0.903"
Set-CellText 40 5 "  -5.62%  "

# Row 41 - WEMIXToken
Set-CellText 41 5 "  +0.77%  "

# Row 42 - PaxDollar
Set-CellText 42 5 "  +0.10%  "

# Row 43 - RenderToken
Set-CellText 43 5 "  +6.95%  "

# Row 44 - FraxShare->Aave
Set-CellText 44 2 "Aave"
Set-CellText 44 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText 44 4 "66.47"
Set-CellText 44 5 "  -2.21%  "

# Row 45 - Aave->FraxShare
Set-CellText 45 2 "FraxShare"
Set-CellText 45 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText 45 4 "5.52"
Set-CellText 45 5 "  +3.33%  "

# Row 46 - MXToken
Set-CellText 46 5 "  -0.02%  "

# Row 47 - RocketPoolETH
Set-CellText 47 4 "1.781.55"
Set-CellText 47 5 "  +0.42%  "

# Row 48 - Quant
Set-CellText 48 4 "87.96"
Set-CellText 48 5 "  +0.16%  "

# Row 49 - Algorand
Set-CellText 49 5 "  +0.96%  "

# Row 50 - Cronos
Set-CellText 50 5 "  +0.30%  "

# Row 51 - EnergySwap
Set-CellText 51 4 "7.60"
Set-CellText 51 5 "  -1.86%  "

